# Standard_Balance_Sheet.xlsx refresh:
#  - insert a new row 1 with a note about the date-header formats the
#    importer understands, merged across A1:D1, styled italic + red
#  - the former row 1 (header labels) becomes row 2; its year labels get
#    an explicit "Annual" suffix and the header style now lives on its
#    own font/xf (no longer shared with the bold xf that used to do the
#    centering+bold combo)
#  - everything else (the account-name rows) simply rides down one row
#    because of the insert, no further edits needed there

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push every existing row down by one, leaving a blank (unstyled) row 1.
$ws.Rows.Item(1).Insert()

# New row 1: explanatory note, italic + red font, merged across A:D.
$ws.Range("A1").Value = "Note: The date header (Row 2) supports: '2023 Annual', '2023 Q1', '2023-01'"
$ws.Range("A1:D1").Merge()
$ws.Range("A1").Font.Italic = $true
$ws.Range("A1").Font.Bold = $false
$ws.Range("A1").Font.Color = 255

# Row 2 (the old row 1, shifted down) keeps the bold+centered header look,
# but the year columns now spell out "Annual".
$ws.Range("B2").Value = "2024 Annual"
$ws.Range("C2").Value = "2023 Annual"
$ws.Range("D2").Value = "2022 Annual"

Write-Host "done"
